# Update "想去人数" (number of people interested) figures that changed
# between the previous gh-pages data pull and this one.
#
# Sheet "展览"   (Exhibitions) : F2 1437->1438, F3 3023->3027, F5 298->340
# Sheet "全部类型" (All types)  : F2 1437->1438, F3 3023->3027, F5 298->340

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1438
    $ws.Range("F3").Value = 3027
    $ws.Range("F5").Value = 340
}
